$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B233").Value = 64979
$ws.Range("E233").Value = 314.41
$ws.Range("F233").Value = 0
$ws.Range("G233").Value = 0
$ws.Range("B234").Value = 48719
$ws.Range("E234").Value = 353.35
$ws.Range("F234").Value = -81
$ws.Range("G234").Value = -23955.75
$ws.Range("B246").Value = 48706
$ws.Range("E246").Value = 39.8
$ws.Range("F246").Value = -144
$ws.Range("G246").Value = -4795.2
$ws.Range("B247").Value = 64973
$ws.Range("E247").Value = 35.4
$ws.Range("F247").Value = 2
$ws.Range("G247").Value = 66.59999999999999
$ws.Range("B292").Value = 55373
$ws.Range("E292").Value = 163.62
$ws.Range("F292").Value = -94
$ws.Range("G292").Value = -13562.32
$ws.Range("B293").Value = 63520
$ws.Range("E293").Value = 153.4
$ws.Range("F293").Value = 69
$ws.Range("G293").Value = 9955.32
$ws.Range("B294").Value = 57802
$ws.Range("E294").Value = 162.71
$ws.Range("F294").Value = -79
$ws.Range("G294").Value = -11334.92
$ws.Range("B295").Value = 63531
$ws.Range("F295").Value = 79
$ws.Range("G295").Value = 11334.92
$ws.Range("B296").Value = 63571
$ws.Range("E296").Value = 152.53
$ws.Range("F296").Value = 0
$ws.Range("G296").Value = 0
$ws.Range("B299").Value = 63510
$ws.Range("E299").Value = 50.66
$ws.Range("F299").Value = 132
$ws.Range("G299").Value = 6288.48
$ws.Range("B300").Value = 55356
$ws.Range("E300").Value = 54.04
$ws.Range("F300").Value = -158
$ws.Range("G300").Value = -7527.12
$ws.Range("B311").Value = 61605
$ws.Range("E311").Value = 133.78
$ws.Range("F311").Value = -13
$ws.Range("G311").Value = -1455.48
$ws.Range("B312").Value = 63563
$ws.Range("E312").Value = 119.04
$ws.Range("F312").Value = 0
$ws.Range("G312").Value = 0
$ws.Range("B420").Value = 58047
$ws.Range("D420").Value = 105.54
$ws.Range("E420").Value = 126.1
$ws.Range("F420").Value = 41
$ws.Range("G420").Value = 4327.14
$ws.Range("B421").Value = 47097
$ws.Range("D421").Value = 112.28
$ws.Range("E421").Value = 134.16
$ws.Range("F421").Value = 15
$ws.Range("G421").Value = 1684.2
$ws.Range("B479").Value = 45718
$ws.Range("E479").Value = 19.38
$ws.Range("F479").Value = -294
$ws.Range("G479").Value = -4768.68
$ws.Range("B480").Value = 64927
$ws.Range("E480").Value = 17.26
$ws.Range("F480").Value = 119
$ws.Range("G480").Value = 1930.18
$ws.Range("B485").Value = 64925
$ws.Range("E485").Value = 13.97
$ws.Range("F485").Value = 129
$ws.Range("G485").Value = 1696.35
$ws.Range("B486").Value = 45709
$ws.Range("E486").Value = 15.69
$ws.Range("F486").Value = -300
$ws.Range("G486").Value = -3945
$ws.Range("B490").Value = 65067
$ws.Range("E490").Value = 15.65
$ws.Range("F490").Value = 172
$ws.Range("G490").Value = 2533.56
$ws.Range("B491").Value = 53595
$ws.Range("E491").Value = 17.61
$ws.Range("F491").Value = -335
$ws.Range("G491").Value = -4934.55
$ws.Range("B600").Value = 64830
$ws.Range("E600").Value = 34.9
$ws.Range("F600").Value = 111
$ws.Range("G600").Value = 3644.13
$ws.Range("B601").Value = 60022
$ws.Range("E601").Value = 37.22
$ws.Range("F601").Value = -113
$ws.Range("G601").Value = -3709.79
$ws.Range("B709").Value = 63150
$ws.Range("D709").Value = 75.68000000000001
$ws.Range("E709").Value = 80.45
$ws.Range("F709").Value = 20
$ws.Range("G709").Value = 1513.6
$ws.Range("B710").Value = 61428
$ws.Range("D710").Value = 69.16
$ws.Range("E710").Value = 73.52
$ws.Range("F710").Value = 1
$ws.Range("G710").Value = 69.16
$ws.Range("B736").Value = 65079
$ws.Range("F736").Value = 21
$ws.Range("G736").Value = 858.27
$ws.Range("B737").Value = 65362
$ws.Range("F737").Value = 37
$ws.Range("G737").Value = 1512.19
